# Update the "想去人数" (interest count) column F with refreshed values for
# each worksheet, as produced by a re-run of the gh-pages data scrape.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - F column updates
$wsExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 4771
    3  = 2748
    5  = 2787
    10 = 751
    11 = 477
    12 = 240
    13 = 410
    14 = 1067
    15 = 300
    16 = 90
    17 = 71
    18 = 536
    19 = 536
    20 = 46
    22 = 661
    23 = 740
    24 = 150
    26 = 515
    27 = 21
    28 = 1659
    29 = 1534
    30 = 338
    32 = 1490
    33 = 183
    34 = 2333
    35 = 391
    37 = 610
    38 = 115
    39 = 60
    41 = 785
    43 = 213
    45 = 494
    46 = 39
    47 = 77
    48 = 109
}
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "演出" (Performance) - F column updates
$wsPerformance = $wb.Worksheets.Item("演出")
$performanceUpdates = @{
    2 = 7
    4 = 90
}
foreach ($row in $performanceUpdates.Keys) {
    $wsPerformance.Cells.Item($row, 6).Value = $performanceUpdates[$row]
}

# Sheet "全部类型" (All types) - F column updates
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 4771
    3  = 2748
    4  = 2787
    7  = 751
    8  = 477
    9  = 240
    10 = 410
    11 = 1067
    12 = 300
    13 = 90
    14 = 71
    15 = 536
    16 = 536
    17 = 46
    18 = 661
    19 = 740
    20 = 150
    21 = 90
    22 = 90
    26 = 515
    27 = 21
    28 = 1659
    29 = 1534
    30 = 338
    33 = 2333
    34 = 391
    39 = 610
    40 = 115
    41 = 60
    43 = 785
    46 = 213
    47 = 494
    48 = 39
    49 = 77
    50 = 109
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
